# Update betting-odds figures on "Sheet1" for the three match rows (3, 4, 5)
# to match the latest FlashScore scrape, per the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: Haras El Hodood vs Smouha ---
$ws.Range("G3").Value = 3.6
$ws.Range("H3").Value = 2.82
$ws.Range("I3").Value = 2.22
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 1.88
$ws.Range("N3").Value = 5.2
$ws.Range("P3").Value = 2.32
$ws.Range("Q3").Value = 2.57
$ws.Range("R3").Value = 1.45
$ws.Range("U3").Value = 2.12
$ws.Range("W3").Value = 7.6
$ws.Range("X3").Value = 17.5
$ws.Range("Y3").Value = 13
$ws.Range("Z3").Value = 55
$ws.Range("AC3").Value = 5.2
$ws.Range("AD3").Value = 5.7
$ws.Range("AE3").Value = 18
$ws.Range("AH3").Value = 5.7
$ws.Range("AL3").Value = 22
$ws.Range("AM3").Value = 40
$ws.Range("AN3").Value = 5.3
$ws.Range("AP3").Value = 35
$ws.Range("AT3").Value = 2.27
$ws.Range("AU3").Value = 7.8
$ws.Range("AX3").Value = 12

# --- Row 4: Haka vs SJK ---
$ws.Range("G4").Value = 3.25
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 3.6
$ws.Range("L4").Value = 2.75
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 23
$ws.Range("AD4").Value = 6.5
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 12
$ws.Range("AK4").Value = 21
$ws.Range("AU4").Value = 7

# --- Row 5: Skalica vs Slovan Bratislava ---
$ws.Range("G5").Value = 6.8
$ws.Range("I5").Value = 1.4
$ws.Range("J5").Value = 5.9
$ws.Range("L5").Value = 1.85
$ws.Range("P5").Value = 4.6
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 1.98
$ws.Range("W5").Value = 19
$ws.Range("Y5").Value = 22
$ws.Range("AA5").Value = 70
$ws.Range("AB5").Value = 60
$ws.Range("AE5").Value = 18.5
$ws.Range("AF5").Value = 75
$ws.Range("AK5").Value = 10.25
$ws.Range("AL5").Value = 11.25
$ws.Range("AN5").Value = 8
$ws.Range("AO5").Value = 35
$ws.Range("AP5").Value = 35
$ws.Range("AS5").Value = 400
$ws.Range("AW5").Value = 3.4
$ws.Range("AX5").Value = 6.3
$ws.Range("AZ5").Value = 16.5
$ws.Range("BA5").Value = 37

$wb.Save()
